$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Jesse")

# --- Row 6: the date was originally typed/stored as the malformed text
# "10/3102017". Fix it to a real date value (10/31/2017).
$ws.Range("A6").Value = 43039

# --- New row 7: a new journal entry.
# Set the plain values first so the formula recalculation (C2 = SUM(B4:B200))
# picks them up.
$ws.Range("A7").Value = 43059
$ws.Range("B7").Value = 300

$longText = "Added data member 'next' and function 'setNext' to Object superclass. `nEdited Object assignment in Rooms.h to account for next data members. `nCreated a test building consisting only of rooms, stairs, and doors that has two floors and a basement.`nAdded BuildingInfo.txt to explain input file formatting.`nCorrected location of stairs in roomID f2r410. `nAdded setSymbol function in Objects.h. `nCompleted stairs symbol assignment in Rooms.h. `nAdded basic menu and probability functions in main. `nChanged definitions from in-line to prototype in header and definition in separate cpp for Objects and Rooms.`nAdded notation sections to headers of Objects.h and Rooms.h. "
$ws.Range("C7").Value = $longText

# Copy the row-6 formatting (borders/fill/number format/wrap) down onto the
# freshly populated row 7.
$ws.Range("A6:C6").Copy()
$ws.Range("A7:C7").PasteSpecial(-4122)  # xlPasteFormats

# Match row height to the new wrapped, multi-line content.
$ws.Rows.Item(7).RowHeight = 270.75

# Move the sheet's active-cell selection to the newly added row, then
# reactivate the workbook's original active sheet (Main) so the overall
# workbook selection state is left the way it started.
$ws.Activate()
$ws.Range("G7").Select()
$wb.Worksheets.Item("Main").Activate()
